# Adding test case for Authoring
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- New row 67: new Authoring test case (VerifyDraftPostTabDisplayForZeroDrafts) ---
# Copy the border/font formatting from the existing last row (65) onto the new
# row 67 cells before writing their values, so the new row keeps the same
# look as the rest of the table. Values are written in the same order the
# target file lists the new shared strings in (Description, JIRA ID, TCID)
# so the newly created shared-string entries line up.
$ws.Cells.Item(65, 1).Copy() | Out-Null
$ws.Cells.Item(67, 3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(67, 3).Value2 = "Verfiy that the Drafts Post tab is not displayed when there are no draft posts"

$ws.Cells.Item(65, 1).Copy() | Out-Null
$ws.Cells.Item(67, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(67, 2).Value2 = "OPQA-1198"

$ws.Cells.Item(65, 1).Copy() | Out-Null
$ws.Cells.Item(67, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(67, 1).Value2 = "VerifyDraftPostTabDisplayForZeroDrafts"

$ws.Cells.Item(65, 1).Copy() | Out-Null
$ws.Cells.Item(67, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(67, 4).Value2 = "Y"

$ws.Cells.Item(65, 1).Copy() | Out-Null
$ws.Cells.Item(67, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(67, 5).Value2 = "PASS"

$excel.CutCopyMode = $false

# --- Row 60: Results PASS -> SKIP ---
$ws.Cells.Item(60, 5).Value2 = "SKIP"

# --- Row 64: Runmode cell format change (s=8 -> s=18), Results PASS -> SKIP ---
$ws.Cells.Item(61, 4).Copy() | Out-Null
$ws.Cells.Item(64, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(64, 5).Value2 = "SKIP"

# --- Row 65: Runmode cell format change (s=1 -> s=18), Results PASS -> SKIP ---
$ws.Cells.Item(61, 4).Copy() | Out-Null
$ws.Cells.Item(65, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(65, 5).Value2 = "SKIP"

# --- Row 66: row height -> 45, Runmode cell format change (s=1 -> s=18), Results PASS -> SKIP ---
$ws.Rows.Item(66).RowHeight = 45
$ws.Cells.Item(61, 4).Copy() | Out-Null
$ws.Cells.Item(66, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(66, 5).Value2 = "SKIP"

$excel.CutCopyMode = $false

# --- sheet view state (active selection / scroll position) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 43
$ws.Range("D66").Select()
